$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.180.64'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '3.891.79'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '483.57'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.92'
$ws.Range("E6").Value = '  -3.03%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.739'
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.183'
$ws.Range("E10").Value = '  +9.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000359'
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.92'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.54'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '4.506.07'
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").Value = '3.918.67'
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.23'
$ws.Range("E16").Value = '  -3.74%  '
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.01'
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").Value = '68.195.08'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '429.88'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.55'
$ws.Range("E22").Value = '  +3.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.82'
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.67'
$ws.Range("E24").Value = '  +2.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.99'
$ws.Range("E25").Value = '  +11.76%  '
$ws.Range("E26").Value = '  +3.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.00'
$ws.Range("E27").Value = '  +1.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.44'
$ws.Range("E28").Value = '  -2.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.66'
$ws.Range("E29").Value = '  -4.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '710.76'
$ws.Range("E30").Value = '  -0.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.53'
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.129'
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("E34").Value = '  +9.81%  '
$ws.Range("D35").Value = '0.0₃0878'
$ws.Range("E35").Value = '  -3.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '40.99'
$ws.Range("E36").Value = '  -2.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '60.91'
$ws.Range("E37").Value = '  +3.19%  '
$ws.Range("E38").Value = '  +6.74%  '
$ws.Range("E39").Value = '  -3.77%  '
$ws.Range("B40").Value = 'Dai'
$ws.Range("C40").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.396'
$ws.Range("E41").Value = '  +14.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.96'
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("E43").Value = '  +2.42%  '
$ws.Range("E44").Value = '  -2.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.142'
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.37'
$ws.Range("E46").Value = '  +4.22%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("E49").Value = '  -3.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '144.31'
$ws.Range("E50").Value = '  -1.20%  '
$ws.Range("E51").Value = '  -1.80%  '
